# Generate Report for Handback
# Stamp fresh handoff/handback datetimes for the newly-processed file
# (1e4c2880-118e-4d5c-8048-b7616a3c6a5b) across the Overview, zh-cn and
# de-de worksheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for row 3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-26 18:50:24"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-26 18:50:20"
$wsZhCn.Range("K3").Value = "2016-08-26 18:50:36"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-26 18:50:24"
$wsDeDe.Range("K3").Value = "2016-08-26 18:50:43"
